$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "HORA F" (end time) value for row 6 - 21:55
$ws.Range("C6").Value = 0.91319444444444453

# Fill in the "DESCANSO" (rest) value for row 6 - 0 (00:00)
$ws.Range("F6").Value = 0

# Update the active selection to F7
$ws.Range("F7").Select()
